# Removed some old text
# The "My Info" slide's contact-info placeholder had a paragraph with
# "OneNote (Take Your Angular App Glamping): " followed by a line break and
# a hyperlinked OneDrive URL. That paragraph's content is removed, leaving
# an empty paragraph (matching the following already-empty paragraph).

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    foreach ($shape in $s.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count
            for ($i = 1; $i -le $paraCount; $i++) {
                $para = $tr.Paragraphs($i, 1)
                if ($para.Text -like "OneNote (Take Your Angular App Glamping)*") {
                    # Delete the paragraph's content (runs + line break) but
                    # keep the paragraph mark itself, so it becomes an empty
                    # paragraph rather than merging with its neighbor.
                    $tr.Characters($para.Start, $para.Length - 1).Delete()
                }
            }
        }
    }
}
